# Insert a new data row before current row 114 ("Macroferia Regional de
# Talca" / Kiwi price sheet). This shifts the existing rows 114-215 down to
# 115-216 and leaves a blank row 114 to be filled in with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(114).Insert()

$ws.Range("A114").Value = 5
$ws.Range("B114").Value = "Macroferia Regional de Talca"
$ws.Range("C114").Value = "Maule"
$ws.Range("D114").Value = 44664
$ws.Range("E114").Value = 7
$ws.Range("F114").Value = "Fruta"
$ws.Range("G114").Value = 100101
$ws.Range("H114").Value = "Berries"
$ws.Range("I114").Value = 100101007
$ws.Range("J114").Value = "Kiwi"
$ws.Range("K114").Value = "Hayward"
$ws.Range("L114").Value = "Primera"
$ws.Range("M114").Value = 200
$ws.Range("N114").Value = 10000
$ws.Range("O114").Value = 10000
$ws.Range("P114").Value = 10000
$ws.Range("Q114").Value = "`$/bandeja 18 kilos"
$ws.Range("R114").Value = "Provincia de Curicó"
$ws.Range("S114").Value = 556
$ws.Range("T114").Value = 18
